$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRange, [string]$text)
    $origStyle = $cellRange.Style
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '43.481.43'
Set-TextValue $ws.Range("E2") '  -1.00%  '
Set-TextValue $ws.Range("D3") '2.220.52'
Set-TextValue $ws.Range("E3") '  +0.84%  '
Set-TextValue $ws.Range("E4") '  -0.03%  '
Set-TextValue $ws.Range("D5") '269.82'
Set-TextValue $ws.Range("D6") '92.27'
Set-TextValue $ws.Range("E6") '  +12.47%  '
Set-TextValue $ws.Range("E7") '  -0.69%  '
Set-TextValue $ws.Range("E8") '  -0.08%  '
Set-TextValue $ws.Range("D9") '0.617'
Set-TextValue $ws.Range("E9") '  +2.79%  '
Set-TextValue $ws.Range("D10") '45.83'
Set-TextValue $ws.Range("E10") '  +5.36%  '
Set-TextValue $ws.Range("D11") '0.0938'
Set-TextValue $ws.Range("E11") '  +1.28%  '
Set-TextValue $ws.Range("D12") '8.21'
Set-TextValue $ws.Range("E12") '  +17.23%  '
Set-TextValue $ws.Range("E13") '  +0.97%  '
Set-TextValue $ws.Range("D14") '2.559.76'
Set-TextValue $ws.Range("E14") '  +0.49%  '
Set-TextValue $ws.Range("D15") '15.01'
Set-TextValue $ws.Range("E15") '  +3.95%  '
Set-TextValue $ws.Range("D16") '2.239.41'
Set-TextValue $ws.Range("E16") '  +2.79%  '
Set-TextValue $ws.Range("D17") '0.799'
Set-TextValue $ws.Range("E17") '  +2.70%  '
Set-TextValue $ws.Range("D18") '43.485.27'
Set-TextValue $ws.Range("E18") '  -0.80%  '
Set-TextValue $ws.Range("D19") '0.0000104'
Set-TextValue $ws.Range("E19") '  +0.31%  '
Set-TextValue $ws.Range("E20") '  +0.50%  '
Set-TextValue $ws.Range("D21") '70.33'
Set-TextValue $ws.Range("E21") '  -1.14%  '
Set-TextValue $ws.Range("E22") '  -1.28%  '
Set-TextValue $ws.Range("D23") '232.17'
Set-TextValue $ws.Range("E23") '  +0.15%  '
Set-TextValue $ws.Range("D24") '9.03'
Set-TextValue $ws.Range("E24") '  -2.05%  '
Set-TextValue $ws.Range("D25") '0.999'
Set-TextValue $ws.Range("D26") '11.32'
Set-TextValue $ws.Range("E26") '  +5.54%  '
Set-TextValue $ws.Range("D27") '2.49'
Set-TextValue $ws.Range("E27") '  +11.00%  '
Set-TextValue $ws.Range("E28") '  +5.11%  '
Set-TextValue $ws.Range("D29") '41.04'
Set-TextValue $ws.Range("E29") '  -0.53%  '
Set-TextValue $ws.Range("E30") '  +2.04%  '
Set-TextValue $ws.Range("D31") '172.50'
Set-TextValue $ws.Range("E31") '  -0.11%  '
Set-TextValue $ws.Range("D32") '0.0918'
Set-TextValue $ws.Range("E32") '  +5.84%  '
Set-TextValue $ws.Range("D33") '20.79'
Set-TextValue $ws.Range("E33") '  +1.29%  '
Set-TextValue $ws.Range("E34") '  +2.82%  '
Set-TextValue $ws.Range("E35") '  +0.43%  '
Set-TextValue $ws.Range("E36") '  -2.86%  '
Set-TextValue $ws.Range("D37") '0.0350'
Set-TextValue $ws.Range("E37") '  -2.30%  '
Set-TextValue $ws.Range("D38") '4.29'
Set-TextValue $ws.Range("E38") '  -4.84%  '
Set-TextValue $ws.Range("D39") '3.60'
Set-TextValue $ws.Range("E39") '  +22.75%  '
Set-TextValue $ws.Range("D40") '12.51'
Set-TextValue $ws.Range("E40") '  -5.69%  '
Set-TextValue $ws.Range("E41") '  +9.53%  '
Set-TextValue $ws.Range("D42") '2.15'
Set-TextValue $ws.Range("E42") '  +2.30%  '
Set-TextValue $ws.Range("D43") '63.16'
Set-TextValue $ws.Range("E43") '  +0.70%  '
Set-TextValue $ws.Range("D44") '5.31'
Set-TextValue $ws.Range("E44") '  -4.01%  '
Set-TextValue $ws.Range("D45") '0.0985'
Set-TextValue $ws.Range("E45") '  -0.07%  '
Set-TextValue $ws.Range("D46") '8.36'
Set-TextValue $ws.Range("E46") '  +0.89%  '
Set-TextValue $ws.Range("D47") '100.05'
Set-TextValue $ws.Range("E47") '  -2.11%  '
Set-TextValue $ws.Range("E48") '  +2.91%  '
Set-TextValue $ws.Range("E49") '  +1.06%  '
Set-TextValue $ws.Range("D50") '0.436'
Set-TextValue $ws.Range("E50") '  -1.37%  '
Set-TextValue $ws.Range("D51") '2.444.97'
Set-TextValue $ws.Range("E51") '  +0.43%  '
